$d = $word.ActiveDocument

# Locate the paragraph containing the "11) Demo - Creating DB objects using SSMS" heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*11) Demo*Creating DB objects*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph ('11) Demo - Creating DB objects using SSMS') not found!"
} else {
    # The blank paragraph that immediately follows the heading is where the
    # new paragraph needs to be added (right after it).
    $emptyPara = $target.Next()

    # Insert a brand-new (blank) paragraph right after that blank paragraph,
    # then fill it in with the new sentence.
    $emptyPara.Range.InsertParagraphAfter() | Out-Null
    $insertedPara = $emptyPara.Next()
    $insertedPara.Range.Text = "Something going on with the video.  It is spinning and spinning.  "

    Write-Host "Inserted new paragraph after the '11) Demo' blank line."
}
